$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear stray cell at AMJ1 (far right of header row) so dimension shrinks
$ws.Range("AMJ1").Clear()

# Add new row of data (row 9)
$ws.Range("A9").Value = "base_mental_health"
$ws.Range("B9").Value = "Mental Health Before"
$ws.Range("C9").Value = "Mental Health Before"
$ws.Range("D9").Value = "Background Variables"
$ws.Range("E9").Value = "Background Overview"

# Update active cell selection to E10
$ws.Range("E10").Select()
